$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 720  # ALC!H32 250000590 -> 720
$ws.Cells.Item(32, 9).Value = 200  # ALC!I32 500000350 -> 200
$ws.Cells.Item(32, 10).Value = 850  # ALC!J32 840 -> 850
$ws.Cells.Item(32, 11).Value = 200  # ALC!K32 500000350 -> 200
$ws.Cells.Item(32, 12).Value = 850  # ALC!L32 840 -> 850
$ws.Cells.Item(32, 13).Value = 126  # ALC!M32 -500000024 -> 126
$ws.Cells.Item(32, 14).Value = -1502  # ALC!N32 -1492 -> -1502

$ws.Cells.Item(49, 8).Value = 4141.6665  # ALC!H49 2215.923 -> 4141.6665
$ws.Cells.Item(49, 9).Value = 900  # ALC!I49 413.4 -> 900
$ws.Cells.Item(49, 10).Value = 4790  # ALC!J49 3342.5 -> 4790
$ws.Cells.Item(49, 11).Value = 2700  # ALC!K49 1240.2 -> 2700
$ws.Cells.Item(49, 12).Value = 14370  # ALC!L49 10027.5 -> 14370
$ws.Cells.Item(49, 13).Value = -2564  # ALC!M49 -1104.2 -> -2564
$ws.Cells.Item(49, 14).Value = -14642  # ALC!N49 -10299.5 -> -14642

$ws.Cells.Item(58, 8).Value = 2224.45  # ALC!H58 2249.389 -> 2224.45
$ws.Cells.Item(58, 10).Value = 4759.3335  # ALC!J58 5547.7144 -> 4759.3335
$ws.Cells.Item(58, 12).Value = 14278.0005  # ALC!L58 16643.1432 -> 14278.0005
$ws.Cells.Item(58, 14).Value = -14578.0005  # ALC!N58 -16943.1432 -> -14578.0005

$ws.Cells.Item(137, 8).Value = 1479.9714  # ALC!H137 1556.2188 -> 1479.9714
$ws.Cells.Item(137, 9).Value = 1194.24  # ALC!I137 1223.1666 -> 1194.24
$ws.Cells.Item(137, 10).Value = 2194.3  # ALC!J137 2555.375 -> 2194.3
$ws.Cells.Item(137, 11).Value = 3582.72  # ALC!K137 3669.4998 -> 3582.72
$ws.Cells.Item(137, 12).Value = 6582.900000000001  # ALC!L137 7666.125 -> 6582.900000000001
$ws.Cells.Item(137, 13).Value = -1032.72  # ALC!M137 -1119.4998 -> -1032.72
$ws.Cells.Item(137, 14).Value = -11682.9  # ALC!N137 -12766.125 -> -11682.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(43, 8).Value = 8233.833000000001  # ARM!H43 6569.625 -> 8233.833000000001
$ws.Cells.Item(43, 9).Value = 2800  # ARM!I43 0 -> 2800
$ws.Cells.Item(43, 10).Value = 9320.6  # ARM!J43 6569.625 -> 9320.6
$ws.Cells.Item(43, 11).Value = 2800  # ARM!K43 0 -> 2800
$ws.Cells.Item(43, 12).Value = 9320.6  # ARM!L43 6569.625 -> 9320.6
$ws.Cells.Item(43, 13).Value = -2487  # ARM!M43 None -> -2487
$ws.Cells.Item(43, 14).Value = -9946.6  # ARM!N43 -7195.625 -> -9946.6

$ws.Cells.Item(54, 8).Value = 8000  # ARM!H54 0 -> 8000
$ws.Cells.Item(54, 10).Value = 8000  # ARM!J54 0 -> 8000
$ws.Cells.Item(54, 12).Value = 8000  # ARM!L54 0 -> 8000
$ws.Cells.Item(54, 14).Value = -9538  # ARM!N54 None -> -9538

$ws.Cells.Item(61, 8).Value = 2055.6316  # ARM!H61 2527.4614 -> 2055.6316
$ws.Cells.Item(61, 9).Value = 2009.8125  # ARM!I61 2595.7 -> 2009.8125
$ws.Cells.Item(61, 11).Value = 2009.8125  # ARM!K61 2595.7 -> 2009.8125
$ws.Cells.Item(61, 13).Value = -1797.8125  # ARM!M61 -2383.7 -> -1797.8125

$ws.Cells.Item(74, 8).Value = 1028.9375  # ARM!H74 802.9474 -> 1028.9375
$ws.Cells.Item(74, 9).Value = 1051.3846  # ARM!I74 830.3111 -> 1051.3846
$ws.Cells.Item(74, 10).Value = 931.6667  # ARM!J74 700.3333 -> 931.6667
$ws.Cells.Item(74, 11).Value = 1051.3846  # ARM!K74 830.3111 -> 1051.3846
$ws.Cells.Item(74, 12).Value = 931.6667  # ARM!L74 700.3333 -> 931.6667
$ws.Cells.Item(74, 13).Value = -177.3846000000001  # ARM!M74 43.68889999999999 -> -177.3846000000001
$ws.Cells.Item(74, 14).Value = -2679.6667  # ARM!N74 -2448.3333 -> -2679.6667

$ws.Cells.Item(77, 8).Value = 1028.9375  # ARM!H77 802.9474 -> 1028.9375
$ws.Cells.Item(77, 9).Value = 1051.3846  # ARM!I77 830.3111 -> 1051.3846
$ws.Cells.Item(77, 10).Value = 931.6667  # ARM!J77 700.3333 -> 931.6667
$ws.Cells.Item(77, 11).Value = 5256.923000000001  # ARM!K77 4151.5555 -> 5256.923000000001
$ws.Cells.Item(77, 12).Value = 4658.3335  # ARM!L77 3501.6665 -> 4658.3335
$ws.Cells.Item(77, 13).Value = -888.9230000000007  # ARM!M77 216.4444999999996 -> -888.9230000000007
$ws.Cells.Item(77, 14).Value = -13394.3335  # ARM!N77 -12237.6665 -> -13394.3335

$ws.Cells.Item(132, 8).Value = 169415.86  # ARM!H132 137521.25 -> 169415.86
$ws.Cells.Item(132, 9).Value = 201099.2  # ARM!I132 173433.38 -> 201099.2
$ws.Cells.Item(132, 10).Value = 10999.2  # ARM!J132 7339.75 -> 10999.2
$ws.Cells.Item(132, 11).Value = 603297.6000000001  # ARM!K132 520300.14 -> 603297.6000000001
$ws.Cells.Item(132, 12).Value = 32997.60000000001  # ARM!L132 22019.25 -> 32997.60000000001
$ws.Cells.Item(132, 13).Value = -600767.6000000001  # ARM!M132 -517770.14 -> -600767.6000000001
$ws.Cells.Item(132, 14).Value = -38057.60000000001  # ARM!N132 -27079.25 -> -38057.60000000001

$ws.Cells.Item(136, 8).Value = 2055.6316  # ARM!H136 2527.4614 -> 2055.6316
$ws.Cells.Item(136, 9).Value = 2009.8125  # ARM!I136 2595.7 -> 2009.8125
$ws.Cells.Item(136, 11).Value = 6029.4375  # ARM!K136 7787.099999999999 -> 6029.4375
$ws.Cells.Item(136, 13).Value = -3479.4375  # ARM!M136 -5237.099999999999 -> -3479.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 276078.28  # BSM!H134 604602.4 -> 276078.28
$ws.Cells.Item(134, 9).Value = 303393.2  # BSM!I134 604602.4 -> 303393.2
$ws.Cells.Item(134, 10).Value = 2929  # BSM!J134 0 -> 2929
$ws.Cells.Item(134, 11).Value = 910179.6000000001  # BSM!K134 1813807.2 -> 910179.6000000001
$ws.Cells.Item(134, 12).Value = 8787  # BSM!L134 0 -> 8787
$ws.Cells.Item(134, 13).Value = -907644.6000000001  # BSM!M134 -1811272.2 -> -907644.6000000001
$ws.Cells.Item(134, 14).Value = -13857  # BSM!N134 None -> -13857

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2751.44  # CRP!H31 2269.6 -> 2751.44
$ws.Cells.Item(31, 9).Value = 2439.2  # CRP!I31 1878.8 -> 2439.2
$ws.Cells.Item(31, 10).Value = 2959.6  # CRP!J31 2660.4 -> 2959.6
$ws.Cells.Item(31, 11).Value = 2439.2  # CRP!K31 1878.8 -> 2439.2
$ws.Cells.Item(31, 12).Value = 2959.6  # CRP!L31 2660.4 -> 2959.6
$ws.Cells.Item(31, 13).Value = -2144.2  # CRP!M31 -1583.8 -> -2144.2
$ws.Cells.Item(31, 14).Value = -3549.6  # CRP!N31 -3250.4 -> -3549.6

$ws.Cells.Item(34, 8).Value = 2751.44  # CRP!H34 2269.6 -> 2751.44
$ws.Cells.Item(34, 9).Value = 2439.2  # CRP!I34 1878.8 -> 2439.2
$ws.Cells.Item(34, 10).Value = 2959.6  # CRP!J34 2660.4 -> 2959.6
$ws.Cells.Item(34, 11).Value = 2439.2  # CRP!K34 1878.8 -> 2439.2
$ws.Cells.Item(34, 12).Value = 2959.6  # CRP!L34 2660.4 -> 2959.6
$ws.Cells.Item(34, 13).Value = -2237.2  # CRP!M34 -1676.8 -> -2237.2
$ws.Cells.Item(34, 14).Value = -3363.6  # CRP!N34 -3064.4 -> -3363.6

$ws.Cells.Item(58, 8).Value = 1658.6666  # CRP!H58 1523.421 -> 1658.6666
$ws.Cells.Item(58, 9).Value = 1658.6666  # CRP!I58 1595.6875 -> 1658.6666
$ws.Cells.Item(58, 10).Value = 0  # CRP!J58 1138 -> 0
$ws.Cells.Item(58, 11).Value = 1658.6666  # CRP!K58 1595.6875 -> 1658.6666
$ws.Cells.Item(58, 12).Value = 0  # CRP!L58 1138 -> 0
$ws.Cells.Item(58, 13).Value = -1455.6666  # CRP!M58 -1392.6875 -> -1455.6666
$ws.Cells.Item(58, 14).ClearContents()  # CRP!N58 remove (was -1544)

$ws.Cells.Item(132, 8).Value = 2800.125  # CRP!H132 2552.6843 -> 2800.125
$ws.Cells.Item(132, 9).Value = 2414.5  # CRP!I132 2206 -> 2414.5
$ws.Cells.Item(132, 11).Value = 7243.5  # CRP!K132 6618 -> 7243.5
$ws.Cells.Item(132, 13).Value = -4713.5  # CRP!M132 -4088 -> -4713.5

$ws.Cells.Item(134, 8).Value = 7589.5264  # CRP!H134 6303.2915 -> 7589.5264
$ws.Cells.Item(134, 9).Value = 7733.3887  # CRP!I134 6775.1904 -> 7733.3887
$ws.Cells.Item(134, 10).Value = 5000  # CRP!J134 3000 -> 5000
$ws.Cells.Item(134, 11).Value = 23200.1661  # CRP!K134 20325.5712 -> 23200.1661
$ws.Cells.Item(134, 12).Value = 15000  # CRP!L134 9000 -> 15000
$ws.Cells.Item(134, 13).Value = -20665.1661  # CRP!M134 -17790.5712 -> -20665.1661
$ws.Cells.Item(134, 14).Value = -20070  # CRP!N134 -14070 -> -20070

$ws.Cells.Item(136, 8).Value = 1658.6666  # CRP!H136 1523.421 -> 1658.6666
$ws.Cells.Item(136, 9).Value = 1658.6666  # CRP!I136 1595.6875 -> 1658.6666
$ws.Cells.Item(136, 10).Value = 0  # CRP!J136 1138 -> 0
$ws.Cells.Item(136, 11).Value = 4975.9998  # CRP!K136 4787.0625 -> 4975.9998
$ws.Cells.Item(136, 12).Value = 0  # CRP!L136 3414 -> 0
$ws.Cells.Item(136, 13).Value = -2425.9998  # CRP!M136 -2237.0625 -> -2425.9998
$ws.Cells.Item(136, 14).ClearContents()  # CRP!N136 remove (was -8514)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 6162  # CUL!H64 6469.952 -> 6162
$ws.Cells.Item(64, 9).Value = 461.25  # CUL!I64 515 -> 461.25
$ws.Cells.Item(64, 10).Value = 7428.8335  # CUL!J64 7462.4443 -> 7428.8335
$ws.Cells.Item(64, 11).Value = 1383.75  # CUL!K64 1545 -> 1383.75
$ws.Cells.Item(64, 12).Value = 22286.5005  # CUL!L64 22387.3329 -> 22286.5005
$ws.Cells.Item(64, 13).Value = -1113.75  # CUL!M64 -1275 -> -1113.75
$ws.Cells.Item(64, 14).Value = -22826.5005  # CUL!N64 -22927.3329 -> -22826.5005

$ws.Cells.Item(67, 8).Value = 6162  # CUL!H67 6469.952 -> 6162
$ws.Cells.Item(67, 9).Value = 461.25  # CUL!I67 515 -> 461.25
$ws.Cells.Item(67, 10).Value = 7428.8335  # CUL!J67 7462.4443 -> 7428.8335
$ws.Cells.Item(67, 11).Value = 1383.75  # CUL!K67 1545 -> 1383.75
$ws.Cells.Item(67, 12).Value = 22286.5005  # CUL!L67 22387.3329 -> 22286.5005
$ws.Cells.Item(67, 13).Value = -447.75  # CUL!M67 -609 -> -447.75
$ws.Cells.Item(67, 14).Value = -24158.5005  # CUL!N67 -24259.3329 -> -24158.5005

$ws.Cells.Item(131, 8).Value = 816.34784  # CUL!H131 843.35486 -> 816.34784
$ws.Cells.Item(131, 9).Value = 246.2  # CUL!I131 270 -> 246.2
$ws.Cells.Item(131, 10).Value = 860.8905999999999  # CUL!J131 862.4666999999999 -> 860.8905999999999
$ws.Cells.Item(131, 11).Value = 738.5999999999999  # CUL!K131 810 -> 738.5999999999999
$ws.Cells.Item(131, 12).Value = 2582.6718  # CUL!L131 2587.4001 -> 2582.6718
$ws.Cells.Item(131, 13).Value = 4301.4  # CUL!M131 4230 -> 4301.4
$ws.Cells.Item(131, 14).Value = -12662.6718  # CUL!N131 -12667.4001 -> -12662.6718

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2004.4375  # GSM!H132 2554.7727 -> 2004.4375
$ws.Cells.Item(132, 9).Value = 1644.4286  # GSM!I132 2068.7896 -> 1644.4286
$ws.Cells.Item(132, 10).Value = 4524.5  # GSM!J132 5632.6665 -> 4524.5
$ws.Cells.Item(132, 11).Value = 4933.2858  # GSM!K132 6206.3688 -> 4933.2858
$ws.Cells.Item(132, 12).Value = 13573.5  # GSM!L132 16897.9995 -> 13573.5
$ws.Cells.Item(132, 13).Value = -2403.2858  # GSM!M132 -3676.3688 -> -2403.2858
$ws.Cells.Item(132, 14).Value = -18633.5  # GSM!N132 -21957.9995 -> -18633.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 889.0417  # LTW!H16 730.63336 -> 889.0417
$ws.Cells.Item(16, 9).Value = 917.45  # LTW!I16 728.11536 -> 917.45
$ws.Cells.Item(16, 11).Value = 917.45  # LTW!K16 728.11536 -> 917.45
$ws.Cells.Item(16, 13).Value = -747.45  # LTW!M16 -558.11536 -> -747.45

$ws.Cells.Item(46, 8).Value = 1506.8  # LTW!H46 1252.8182 -> 1506.8
$ws.Cells.Item(46, 9).Value = 1200  # LTW!I46 983.5714 -> 1200
$ws.Cells.Item(46, 10).Value = 1857.4286  # LTW!J46 1724 -> 1857.4286
$ws.Cells.Item(46, 11).Value = 1200  # LTW!K46 983.5714 -> 1200
$ws.Cells.Item(46, 12).Value = 1857.4286  # LTW!L46 1724 -> 1857.4286
$ws.Cells.Item(46, 13).Value = -1012  # LTW!M46 -795.5714 -> -1012
$ws.Cells.Item(46, 14).Value = -2233.4286  # LTW!N46 -2100 -> -2233.4286

$ws.Cells.Item(132, 8).Value = 3306.5334  # LTW!H132 2905.111 -> 3306.5334
$ws.Cells.Item(132, 9).Value = 2111.111  # LTW!I132 1807.8334 -> 2111.111
$ws.Cells.Item(132, 11).Value = 6333.333  # LTW!K132 5423.5002 -> 6333.333
$ws.Cells.Item(132, 13).Value = -3803.333  # LTW!M132 -2893.5002 -> -3803.333

$ws.Cells.Item(136, 8).Value = 1821.5416  # LTW!H136 1810.85 -> 1821.5416
$ws.Cells.Item(136, 9).Value = 1585.85  # LTW!I136 1643 -> 1585.85
$ws.Cells.Item(136, 10).Value = 3000  # LTW!J136 5000 -> 3000
$ws.Cells.Item(136, 11).Value = 4757.549999999999  # LTW!K136 4929 -> 4757.549999999999
$ws.Cells.Item(136, 12).Value = 9000  # LTW!L136 15000 -> 9000
$ws.Cells.Item(136, 13).Value = -2207.549999999999  # LTW!M136 -2379 -> -2207.549999999999
$ws.Cells.Item(136, 14).Value = -14100  # LTW!N136 -20100 -> -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 15149.8  # WVR!H39 16824.834 -> 15149.8
$ws.Cells.Item(39, 9).Value = 2800  # WVR!I39 11000 -> 2800
$ws.Cells.Item(39, 10).Value = 18237.25  # WVR!J39 17989.8 -> 18237.25
$ws.Cells.Item(39, 11).Value = 2800  # WVR!K39 11000 -> 2800
$ws.Cells.Item(39, 12).Value = 18237.25  # WVR!L39 17989.8 -> 18237.25
$ws.Cells.Item(39, 13).Value = -2387  # WVR!M39 -10587 -> -2387
$ws.Cells.Item(39, 14).Value = -19063.25  # WVR!N39 -18815.8 -> -19063.25

$ws.Cells.Item(41, 8).Value = 5996.5  # WVR!H41 6134.8335 -> 5996.5
$ws.Cells.Item(41, 10).Value = 5996.5  # WVR!J41 6134.8335 -> 5996.5
$ws.Cells.Item(41, 12).Value = 5996.5  # WVR!L41 6134.8335 -> 5996.5
$ws.Cells.Item(41, 14).Value = -6776.5  # WVR!N41 -6914.8335 -> -6776.5

$ws.Cells.Item(45, 8).Value = 6429.5454  # WVR!H45 11526 -> 6429.5454
$ws.Cells.Item(45, 9).Value = 8569  # WVR!I45 0 -> 8569
$ws.Cells.Item(45, 10).Value = 6215.6  # WVR!J45 11526 -> 6215.6
$ws.Cells.Item(45, 11).Value = 8569  # WVR!K45 0 -> 8569
$ws.Cells.Item(45, 12).Value = 6215.6  # WVR!L45 11526 -> 6215.6
$ws.Cells.Item(45, 13).Value = -8078  # WVR!M45 None -> -8078
$ws.Cells.Item(45, 14).Value = -7197.6  # WVR!N45 -12508 -> -7197.6

$ws.Cells.Item(74, 8).Value = 10313  # WVR!H74 6026.5 -> 10313
$ws.Cells.Item(74, 9).Value = 0  # WVR!I74 3500 -> 0
$ws.Cells.Item(74, 10).Value = 10313  # WVR!J74 6868.6665 -> 10313
$ws.Cells.Item(74, 11).Value = 0  # WVR!K74 3500 -> 0
$ws.Cells.Item(74, 12).Value = 10313  # WVR!L74 6868.6665 -> 10313
$ws.Cells.Item(74, 13).ClearContents()  # WVR!M74 remove (was -2564)
$ws.Cells.Item(74, 14).Value = -12185  # WVR!N74 -8740.666499999999 -> -12185

$ws.Cells.Item(77, 8).Value = 10313  # WVR!H77 6026.5 -> 10313
$ws.Cells.Item(77, 9).Value = 0  # WVR!I77 3500 -> 0
$ws.Cells.Item(77, 10).Value = 10313  # WVR!J77 6868.6665 -> 10313
$ws.Cells.Item(77, 11).Value = 0  # WVR!K77 10500 -> 0
$ws.Cells.Item(77, 12).Value = 30939  # WVR!L77 20605.9995 -> 30939
$ws.Cells.Item(77, 13).ClearContents()  # WVR!M77 remove (was -5820)
$ws.Cells.Item(77, 14).Value = -40299  # WVR!N77 -29965.9995 -> -40299

$ws.Cells.Item(132, 8).Value = 3955.1853  # WVR!H132 4729.909 -> 3955.1853
$ws.Cells.Item(132, 9).Value = 3722.4092  # WVR!I132 4656.5293 -> 3722.4092
$ws.Cells.Item(132, 11).Value = 11167.2276  # WVR!K132 13969.5879 -> 11167.2276
$ws.Cells.Item(132, 13).Value = -8637.2276  # WVR!M132 -11439.5879 -> -8637.2276

$ws.Cells.Item(136, 8).Value = 1416.3715  # WVR!H136 1467.2941 -> 1416.3715
$ws.Cells.Item(136, 9).Value = 1223.3235  # WVR!I136 1297.125 -> 1223.3235
$ws.Cells.Item(136, 10).Value = 7980  # WVR!J136 4190 -> 7980
$ws.Cells.Item(136, 11).Value = 3669.9705  # WVR!K136 3891.375 -> 3669.9705
$ws.Cells.Item(136, 12).Value = 23940  # WVR!L136 12570 -> 23940
$ws.Cells.Item(136, 13).Value = -1119.9705  # WVR!M136 -1341.375 -> -1119.9705
$ws.Cells.Item(136, 14).Value = -29040  # WVR!N136 -17670 -> -29040
